$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.257450056250921
$ws.Range("D2").Value = 2.795147031779273
$ws.Range("E2").Value = 10.83537105113994
$ws.Range("F2").Value = 21.60762377134271
$ws.Range("G2").Value = 3.570010303265964
$ws.Range("I2").Value = 17.4944546962406
$ws.Range("M2").Value = 16.86830180225451
$ws.Range("N2").Value = 16.87520808968837
$ws.Range("O2").Value = 18.54145570389145
$ws.Range("B3").Value = 7.059030704364075
$ws.Range("D3").Value = 2.802104594502937
$ws.Range("E3").Value = 10.95337658070012
$ws.Range("F3").Value = 21.1431742536105
$ws.Range("G3").Value = 3.572982899121306
$ws.Range("I3").Value = 17.54173133322975
$ws.Range("M3").Value = 16.23353069052102
$ws.Range("N3").Value = 16.82681263313744
$ws.Range("O3").Value = 18.26382946549023
$ws.Range("B4").Value = 6.934755501531755
$ws.Range("D4").Value = 2.806522493536974
$ws.Range("E4").Value = 11.03054916534537
$ws.Range("F4").Value = 20.8606604217961
$ws.Range("G4").Value = 3.574904617068676
$ws.Range("I4").Value = 17.57457754136478
$ws.Range("M4").Value = 15.8330711076163
$ws.Range("N4").Value = 16.79910849912733
$ws.Range("O4").Value = 18.09766987823232
$ws.Range("B5").Value = 6.883567754431458
$ws.Range("D5").Value = 2.808359771940897
$ws.Range("E5").Value = 11.06317610185952
$ws.Range("F5").Value = 20.74638941896509
$ws.Range("G5").Value = 3.575712093469121
$ws.Range("I5").Value = 17.58892312291418
$ws.Range("M5").Value = 15.66743349923294
$ws.Range("N5").Value = 16.78833308416154
$ws.Range("O5").Value = 18.03113212615506
$ws.Range("B6").Value = 6.875037316149321
$ws.Range("D6").Value = 2.808667089689541
$ws.Range("E6").Value = 11.068664742792
$ws.Range("F6").Value = 20.72747197739895
$ws.Range("G6").Value = 3.57584764819337
$ws.Range("I6").Value = 17.59136321535756
$ws.Range("M6").Value = 15.63978938026638
$ws.Range("N6").Value = 16.78657513773711
$ws.Range("O6").Value = 18.02015708596705
$ws.Range("B7").Value = 6.934067274620707
$ws.Range("D7").Value = 2.806547121795373
$ws.Range("E7").Value = 11.03098442179862
$ws.Range("F7").Value = 20.85911561051615
$ws.Range("G7").Value = 3.574915408227825
$ws.Range("I7").Value = 17.574767121516
$ws.Range("M7").Value = 15.83084682587201
$ws.Range("N7").Value = 16.79896108488262
$ws.Range("O7").Value = 18.09676765616009
$ws.Range("B8").Value = 7.18958320541039
$ws.Range("D8").Value = 2.797515865353799
$ws.Range("E8").Value = 10.87507574357923
$ws.Range("F8").Value = 21.44704087078534
$ws.Range("G8").Value = 3.571015271512572
$ws.Range("I8").Value = 17.50996348332879
$ws.Range("M8").Value = 16.65179442033955
$ws.Range("N8").Value = 16.85810924070712
$ws.Range("O8").Value = 18.44489239560949
$ws.Range("B9").Value = 7.668401922233094
$ws.Range("D9").Value = 2.780951803854453
$ws.Range("E9").Value = 10.60711962677258
$ws.Range("F9").Value = 22.61300671985516
$ws.Range("G9").Value = 3.56412909235062
$ws.Range("I9").Value = 17.4131586186219
$ws.Range("M9").Value = 18.16720803414243
$ws.Range("N9").Value = 16.98966247061546
$ws.Range("O9").Value = 19.15760072497103
$ws.Range("B10").Value = 8.003153337858365
$ws.Range("D10").Value = 2.769464515799487
$ws.Range("E10").Value = 10.43379784381119
$ws.Range("F10").Value = 23.46714987781533
$ws.Range("G10").Value = 3.559528801166555
$ws.Range("I10").Value = 17.36045864629023
$ws.Range("M10").Value = 19.21196392321212
$ws.Range("N10").Value = 17.09525872400825
$ws.Range("O10").Value = 19.69396108870488
$ws.Range("B11").Value = 8.151103436173326
$ws.Range("D11").Value = 2.76438330812486
$ws.Range("E11").Value = 10.36018232163866
$ws.Range("F11").Value = 23.85310730628982
$ws.Range("G11").Value = 3.557534493357244
$ws.Range("I11").Value = 17.34047431257432
$ws.Range("M11").Value = 19.67046519814816
$ws.Range("N11").Value = 17.14511686432284
$ws.Range("O11").Value = 19.93955160537302
$ws.Range("B12").Value = 8.206458867138199
$ws.Range("D12").Value = 2.762479688100261
$ws.Range("E12").Value = 10.3330688852239
$ws.Range("F12").Value = 23.99872502443833
$ws.Range("G12").Value = 3.556793358513465
$ws.Range("I12").Value = 17.33347929788431
$ws.Range("M12").Value = 19.84154657403322
$ws.Range("N12").Value = 17.1642479815489
$ws.Range("O12").Value = 20.03268373712908
$ws.Range("B13").Value = 8.194567592293927
$ws.Range("D13").Value = 2.762888757821873
$ws.Range("E13").Value = 10.3388741228056
$ws.Range("F13").Value = 23.96738996627508
$ws.Range("G13").Value = 3.556952350860165
$ws.Range("I13").Value = 17.33496035002617
$ws.Range("M13").Value = 19.80481616649417
$ws.Range("N13").Value = 17.16011677223656
$ws.Range("O13").Value = 20.01262170392893
$ws.Range("B14").Value = 8.15567124517276
$ws.Range("D14").Value = 2.764226286238585
$ws.Range("E14").Value = 10.35793632188978
$ws.Range("F14").Value = 23.86509900212709
$ws.Range("G14").Value = 3.557473238306331
$ws.Range("I14").Value = 17.33988735878107
$ws.Range("M14").Value = 19.68459182969203
$ws.Range("N14").Value = 17.14668582829338
$ws.Range("O14").Value = 19.94721153851029
$ws.Range("B15").Value = 8.13175751862167
$ws.Range("D15").Value = 2.765048226512727
$ws.Range("E15").Value = 10.36971219729183
$ws.Range("F15").Value = 23.80236837279156
$ws.Range("G15").Value = 3.557794126327211
$ws.Range("I15").Value = 17.3429798327976
$ws.Range("M15").Value = 19.610616134693
$ws.Range("N15").Value = 17.1384913319563
$ws.Range("O15").Value = 19.90716030994328
$ws.Range("B16").Value = 7.993393162652946
$ws.Range("D16").Value = 2.769799469810113
$ws.Range("E16").Value = 10.43871496335043
$ws.Range("F16").Value = 23.44186028216751
$ws.Range("G16").Value = 3.559661106681892
$ws.Range("I16").Value = 17.36184485828815
$ws.Range("M16").Value = 19.18165072167082
$ws.Range("N16").Value = 17.09203609852296
$ws.Range("O16").Value = 19.67793493632065
$ws.Range("B17").Value = 7.907366898113716
$ws.Range("D17").Value = 2.772751018122636
$ws.Range("E17").Value = 10.48239308100264
$ws.Range("F17").Value = 23.2199164620499
$ws.Range("G17").Value = 3.560831579768865
$ws.Range("I17").Value = 17.37443898633338
$ws.Range("M17").Value = 18.91410073208852
$ws.Range("N17").Value = 17.06399615809472
$ws.Range("O17").Value = 19.53765171397449
$ws.Range("B18").Value = 7.857481842068885
$ws.Range("D18").Value = 2.774462278345906
$ws.Range("E18").Value = 10.50800709405194
$ws.Range("F18").Value = 23.09202497793401
$ws.Range("G18").Value = 3.561514070512552
$ws.Range("I18").Value = 17.38205834436223
$ws.Range("M18").Value = 18.75864191632025
$ws.Range("N18").Value = 17.04804048905712
$ws.Range("O18").Value = 19.45712057448493
$ws.Range("B19").Value = 7.840523575984578
$ws.Range("D19").Value = 2.775044026248953
$ws.Range("E19").Value = 10.51676364554439
$ws.Range("F19").Value = 23.04868782345994
$ws.Range("G19").Value = 3.561746744022565
$ws.Range("I19").Value = 17.38470265713886
$ws.Range("M19").Value = 18.70574072352123
$ws.Range("N19").Value = 17.0426680755098
$ws.Range("O19").Value = 19.42988392421898
$ws.Range("B20").Value = 7.916566817258391
$ws.Range("D20").Value = 2.772435414228762
$ws.Range("E20").Value = 10.47769252186443
$ws.Range("F20").Value = 23.24356831884245
$ws.Range("G20").Value = 3.560706022555011
$ws.Range("I20").Value = 17.37305945892282
$ws.Range("M20").Value = 18.94274550917121
$ws.Range("N20").Value = 17.06696331871938
$ws.Range("O20").Value = 19.55256966431277
$ws.Range("B21").Value = 8.167114585682809
$ws.Range("D21").Value = 2.763832866663671
$ws.Range("E21").Value = 10.35231648794679
$ws.Range("F21").Value = 23.89516013168898
$ws.Range("G21").Value = 3.55731985993268
$ws.Range("I21").Value = 17.33842464508706
$ws.Range("M21").Value = 19.7199746192274
$ws.Range("N21").Value = 17.15062410160538
$ws.Range("O21").Value = 19.96642126165341
$ws.Range("B22").Value = 8.326937589751923
$ws.Range("D22").Value = 2.758330116819894
$ws.Range("E22").Value = 10.27482990540175
$ws.Range("F22").Value = 24.3178204984444
$ws.Range("G22").Value = 3.555188758095162
$ws.Range("I22").Value = 17.31912586320057
$ws.Range("M22").Value = 20.21306628517496
$ws.Range("N22").Value = 17.20675851358769
$ws.Range("O22").Value = 20.23762529568764
$ws.Range("B23").Value = 8.242010627431778
$ws.Range("D23").Value = 2.761256182693748
$ws.Range("E23").Value = 10.31577457722704
$ws.Range("F23").Value = 24.09258199600781
$ws.Range("G23").Value = 3.556318695623742
$ws.Range("I23").Value = 17.3291210227423
$ws.Range("M23").Value = 19.95129409860804
$ws.Range("N23").Value = 17.17666890286031
$ws.Range("O23").Value = 20.09284346498259
$ws.Range("B24").Value = 7.912408861456356
$ws.Range("D24").Value = 2.77257805393569
$ws.Range("E24").Value = 10.47981607905254
$ws.Range("F24").Value = 23.23287621457531
$ws.Range("G24").Value = 3.560762757155453
$ws.Range("I24").Value = 17.3736819632107
$ws.Range("M24").Value = 18.92980030580888
$ws.Range("N24").Value = 17.06562135102925
$ws.Range("O24").Value = 19.54582487908646
$ws.Range("B25").Value = 7.541628136734456
$ws.Range("D25").Value = 2.785311818096341
$ws.Range("E25").Value = 10.6755129957111
$ws.Range("F25").Value = 22.29728869040883
$ws.Range("G25").Value = 3.565910985471263
$ws.Range("I25").Value = 17.43610868257541
$ws.Range("M25").Value = 17.76858941429932
$ws.Range("N25").Value = 16.95246244097758
$ws.Range("O25").Value = 18.96215410841775
